# Write reporting log files: update harvest rows 2-6 and fix species casing on row 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..6) {
    $ws.Range("B$row").Value = "Yes"
    $ws.Range("F$row").Value = "Na"
    $ws.Range("J$row").Value = 1
}

$ws.Range("F7").Value = "Fisher"
